$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to stay a text cell (the sheet stores
# purely numeric-looking values such as quantities/prices as text strings, even
# though the cell's number format looks numeric), without disturbing the
# cell's existing number-format / style.
function Set-TextValue {
    param($range, $value)
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# The first item (BRUFEN 100MG/5ML SYRUP 150ML) was removed from the list.
# Every remaining item shifts up by one row (row 8 -> row 7, 9 -> 8, 10 -> 9, 11 -> 10),
# and a brand new item (قطن 50جم) is appended as the new last row (row 11).
# Columns: C = item name, H = "الرصيد الحالي", L = "حد الطلب",
#          N = "السعر", P = "سعر البيع", Q = "عدد التعااملات"
# C, H and Q hold values that are never pure numbers (names or "x:y" counters),
# so a plain assignment keeps them as text automatically. L, N and P hold plain
# numeric-looking text ("1", "40.00", ...), which Excel would otherwise coerce
# into real numbers, so those use the text-forcing helper.

$ws.Range("C7").Value = "COENZYME Q10 30 MG 20 CAPS."
$ws.Range("H7").Value = "0:0"
Set-TextValue $ws.Range("L7") "1"
Set-TextValue $ws.Range("N7") "40.00"
Set-TextValue $ws.Range("P7") "40.0000"
$ws.Range("Q7").Value = "1:0"

$ws.Range("C8").Value = "DOXIRAZOL 60 MG 14 CAPS."
$ws.Range("H8").Value = "0:0"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "101.00"
Set-TextValue $ws.Range("P8") "101.0000"
$ws.Range("Q8").Value = "1:0"

$ws.Range("C9").Value = "MIDODRINE 2.5MG 20 TAB"
$ws.Range("H9").Value = "0:1"
Set-TextValue $ws.Range("L9") "1"
Set-TextValue $ws.Range("N9") "50.00"
Set-TextValue $ws.Range("P9") "25.0000"
$ws.Range("Q9").Value = "0:1"

$ws.Range("C10").Value = "زيت جونسون صغير"
$ws.Range("H10").Value = "3:0"
Set-TextValue $ws.Range("L10") "0"
Set-TextValue $ws.Range("N10") "65.00"
Set-TextValue $ws.Range("P10") "65.0000"
$ws.Range("Q10").Value = "1:0"

$ws.Range("C11").Value = "قطن 50جم"
$ws.Range("H11").Value = "15:0"
Set-TextValue $ws.Range("L11") "0"
Set-TextValue $ws.Range("N11") "10.00"
Set-TextValue $ws.Range("P11") "10.0000"
$ws.Range("Q11").Value = "1:0"

# Update the total in P12
$ws.Range("P12").Value = 241

# Update the footer timestamp to reflect the new export time
$ws.Range("A13").Value = "Sunday, 22 June, 2025 9:54 AM"
